$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnas"
$ws.Range("C2").Value = "Gcgr"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 84.76851766666665
$ws.Range("H2").Value = 254.305553
$ws.Range("I2").Value = 0.2571740874301185
$ws.Range("J2").Value = 0.2571740874301185
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8254096666666667
$ws.Range("N2").Value = 2.476229
$ws.Range("O2").Value = 0.9632797185740577
$ws.Range("P2").Value = 0.9632797185740578
$ws.Range("Q2").Value = 69.96875391107078
$ws.Range("R2").Value = 629.718785199637
$ws.Range("S2").Value = 0.2477305825642246
$ws.Range("T2").Value = 0.2477305825642247

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnas"
$ws.Range("C3").Value = "Gcgr"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 84.76851766666665
$ws.Range("H3").Value = 254.305553
$ws.Range("I3").Value = 0.2571740874301185
$ws.Range("J3").Value = 0.2571740874301185
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.03146466666666667
$ws.Range("N3").Value = 0.094394
$ws.Range("O3").Value = 0.03672028142594227
$ws.Range("P3").Value = 0.03672028142594227
$ws.Range("Q3").Value = 2.667213152209111
$ws.Range("R3").Value = 24.004918369882
$ws.Range("S3").Value = 0.009443504865893832
$ws.Range("T3").Value = 0.009443504865893836

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gnas"
$ws.Range("C4").Value = "Gcgr"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 133.6830266666667
$ws.Range("H4").Value = 401.04908
$ws.Range("I4").Value = 0.4055728628296552
$ws.Range("J4").Value = 0.4055728628296552
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8254096666666667
$ws.Range("N4").Value = 2.476229
$ws.Range("O4").Value = 0.9632797185740577
$ws.Range("P4").Value = 0.9632797185740578
$ws.Range("Q4").Value = 110.3432624799245
$ws.Range("R4").Value = 993.08936231932
$ws.Range("S4").Value = 0.3906801131678251
$ws.Range("T4").Value = 0.3906801131678252

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnas"
$ws.Range("C5").Value = "Gcgr"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 133.6830266666667
$ws.Range("H5").Value = 401.04908
$ws.Range("I5").Value = 0.4055728628296552
$ws.Range("J5").Value = 0.4055728628296552
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03146466666666667
$ws.Range("N5").Value = 0.094394
$ws.Range("O5").Value = 0.03672028142594227
$ws.Range("P5").Value = 0.03672028142594227
$ws.Range("Q5").Value = 4.206291873057778
$ws.Range("R5").Value = 37.85662685752001
$ws.Range("S5").Value = 0.01489274966183002
$ws.Range("T5").Value = 0.01489274966183002

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Gnas"
$ws.Range("C6").Value = "Gcgr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 111.16377
$ws.Range("H6").Value = 333.49131
$ws.Range("I6").Value = 0.3372530497402263
$ws.Range("J6").Value = 0.3372530497402264
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8254096666666667
$ws.Range("N6").Value = 2.476229
$ws.Range("O6").Value = 0.9632797185740577
$ws.Range("P6").Value = 0.9632797185740578
$ws.Range("Q6").Value = 91.75565034111
$ws.Range("R6").Value = 825.80085306999
$ws.Range("S6").Value = 0.3248690228420079
$ws.Range("T6").Value = 0.324869022842008

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Gnas"
$ws.Range("C7").Value = "Gcgr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 111.16377
$ws.Range("H7").Value = 333.49131
$ws.Range("I7").Value = 0.3372530497402263
$ws.Range("J7").Value = 0.3372530497402264
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.03146466666666667
$ws.Range("N7").Value = 0.094394
$ws.Range("O7").Value = 0.03672028142594227
$ws.Range("P7").Value = 0.03672028142594227
$ws.Range("Q7").Value = 3.49773096846
$ws.Range("R7").Value = 31.47957871614
$ws.Range("S7").Value = 0.01238402689821842
$ws.Range("T7").Value = 0.01238402689821842
